$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 corresponds to cm008. Mark it as linked ("link_it" column C) and set its
# topic (column D) to the new "Discrete random variables in R" session.
$ws.Range("C9").Value = $true
$ws.Range("D9").Value = "Discrete random variables in R"

# Update the active selection to match the authored state.
$ws.Range("D10").Select()
